{"js": "// Replace each three-digit \u00f7 one-digit division prompt in the worksheet\n// table with its newly generated counterpart, per the commit's mapping.\nconst replacements = [\n  [\"165\u00f75=\", \"429\u00f76=\"],\n  [\"111\u00f73=\", \"890\u00f79=\"],\n  [\"155\u00f79=\", \"812\u00f73=\"],\n  [\"848\u00f79=\", \"960\u00f77=\"],\n  [\"514\u00f72=\", \"877\u00f78=\"],\n  [\"180\u00f75=\", \"301\u00f72=\"],\n  [\"123\u00f77=\", \"108\u00f79=\"],\n  [\"672\u00f78=\", \"782\u00f73=\"],\n  [\"760\u00f78=\", \"688\u00f79=\"],\n  [\"467\u00f78=\", \"581\u00f74=\"],\n  [\"411\u00f77=\", \"328\u00f77=\"],\n  [\"979\u00f78=\", \"674\u00f76=\"],\n  [\"760\u00f76=\", \"693\u00f76=\"],\n  [\"816\u00f74=\", \"449\u00f77=\"],\n  [\"507\u00f76=\", \"889\u00f77=\"],\n  [\"816\u00f77=\", \"581\u00f77=\"],\n  [\"958\u00f79=\", \"975\u00f79=\"],\n  [\"335\u00f76=\", \"509\u00f73=\"],\n  [\"341\u00f75=\", \"335\u00f72=\"],\n  [\"144\u00f75=\", \"972\u00f75=\"],\n  [\"994\u00f79=\", \"830\u00f72=\"],\n  [\"188\u00f76=\", \"830\u00f74=\"],\n  [\"647\u00f72=\", \"751\u00f74=\"],\n  [\"134\u00f72=\", \"494\u00f77=\"],\n  [\"372\u00f77=\", \"915\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit \u00f7 one-digit division prompt in the worksheet\n# table with its newly generated counterpart, per the commit's mapping.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"165\u00f75=\", \"429\u00f76=\"),\n    @(\"111\u00f73=\", \"890\u00f79=\"),\n    @(\"155\u00f79=\", \"812\u00f73=\"),\n    @(\"848\u00f79=\", \"960\u00f77=\"),\n    @(\"514\u00f72=\", \"877\u00f78=\"),\n    @(\"180\u00f75=\", \"301\u00f72=\"),\n    @(\"123\u00f77=\", \"108\u00f79=\"),\n    @(\"672\u00f78=\", \"782\u00f73=\"),\n    @(\"760\u00f78=\", \"688\u00f79=\"),\n    @(\"467\u00f78=\", \"581\u00f74=\"),\n    @(\"411\u00f77=\", \"328\u00f77=\"),\n    @(\"979\u00f78=\", \"674\u00f76=\"),\n    @(\"760\u00f76=\", \"693\u00f76=\"),\n    @(\"816\u00f74=\", \"449\u00f77=\"),\n    @(\"507\u00f76=\", \"889\u00f77=\"),\n    @(\"816\u00f77=\", \"581\u00f77=\"),\n    @(\"958\u00f79=\", \"975\u00f79=\"),\n    @(\"335\u00f76=\", \"509\u00f73=\"),\n    @(\"341\u00f75=\", \"335\u00f72=\"),\n    @(\"144\u00f75=\", \"972\u00f75=\"),\n    @(\"994\u00f79=\", \"830\u00f72=\"),\n    @(\"188\u00f76=\", \"830\u00f74=\"),\n    @(\"647\u00f72=\", \"751\u00f74=\"),\n    @(\"134\u00f72=\", \"494\u00f77=\"),\n    @(\"372\u00f77=\", \"915\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
